$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the old "verlof" (E) and "recup" (F) boolean columns ---
# (rightmost first so column indices stay valid while deleting)
$ws.Columns.Item(6).Delete()
$ws.Columns.Item(5).Delete()

# --- 2. Insert a new row for the missing day (2021-10-04, serial 44473) ---
$ws.Rows.Item(5).Insert()

# --- 3. Give the whole D data column (rows 2-10) the same number format / style
#        as column C (date-time), matching the rest of the sheet ---
$ws.Range("C2").Copy()
$ws.Range("D2:D10").PasteSpecial(-4122)

# --- 4. Header: "wacht" -> "type" ---
$ws.Range("D1").Value = "type"

# --- 5. Fill in the per-day "type" values (single-letter codes) ---
$ws.Range("D4").Value = "W"

$ws.Range("A5").Value = 44473
$ws.Range("D5").Value = "R"

$ws.Range("D8").Value = "R"
$ws.Range("D9").Value = "V"
$ws.Range("D10").Value = "W"

# --- 6. Column D width, matching column C ---
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# --- 7. Update the active selection ---
[void]$ws.Range("E6").Select()
